# Add an "Integrity" column (G) that counts how many of the C:F wire-reading
# cells are populated (numeric) for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell G1 — reuse the same header style (bold/centered/bordered) as
# the existing header cells by copying formatting from F1, then set the text.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Integrity"

# Data rows 2..376 — count populated cells among C:F and write the count to G.
$lastRow = 376
for ($r = 2; $r -le $lastRow; $r++) {
    $count = 0
    for ($c = 3; $c -le 6; $c++) {
        $val = $ws.Cells.Item($r, $c).Value2
        if ($val -ne "") {
            $count = $count + 1
        }
    }
    $ws.Cells.Item($r, 7).Value = $count
}
